# Regenerate the "Report for Archive": update the handoff status text and
# resize the status columns to match the shorter replacement text.

$wb = $excel.ActiveWorkbook

$oldText = "Ready for handoff"
$newText = "In Translation"
# Closest column width (in Excel "characters" units) this runtime can commit
# that matches the narrower, autofit-style width used for the new, shorter
# status text.
$newColumnWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    $rowOffset = $used.Row
    $colOffset = $used.Column

    for ($r = 0; $r -lt $rows; $r++) {
        for ($c = 0; $c -lt $cols; $c++) {
            $cell = $ws.Cells.Item($rowOffset + $r, $colOffset + $c)
            # Note: keep the string literal on the left of -eq. If a boolean
            # cell value (e.g. True/False) ends up on the left instead, the
            # right-hand side gets coerced to a bool and matches everything.
            if ($oldText -eq $cell.Value2) {
                $cell.Value2 = $newText
                $cell.EntireColumn.ColumnWidth = $newColumnWidth
            }
        }
    }
}
